$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Suite1 sheet: flip Run Flag (col B) from N to Y for rows 4,5,7,8,9,10
# (row 6 intentionally stays N)
# ---------------------------------------------------------------------------
$suite1 = $wb.Worksheets.Item("Suite1")
$suite1.Range("B4").Value = "Y"
$suite1.Range("B5").Value = "Y"
$suite1.Range("B7").Value = "Y"
$suite1.Range("B8").Value = "Y"
$suite1.Range("B9").Value = "Y"
$suite1.Range("B10").Value = "Y"

# ---------------------------------------------------------------------------
# Suite2 sheet: flip Run Flag (col B) from N to Y for rows 4 through 20
# ---------------------------------------------------------------------------
$suite2 = $wb.Worksheets.Item("Suite2")
$suite2.Range("B4").Value = "Y"
$suite2.Range("B5").Value = "Y"
$suite2.Range("B6").Value = "Y"
$suite2.Range("B7").Value = "Y"
$suite2.Range("B8").Value = "Y"
$suite2.Range("B9").Value = "Y"
$suite2.Range("B10").Value = "Y"
$suite2.Range("B11").Value = "Y"
$suite2.Range("B12").Value = "Y"
$suite2.Range("B13").Value = "Y"
$suite2.Range("B14").Value = "Y"
$suite2.Range("B15").Value = "Y"
$suite2.Range("B16").Value = "Y"
$suite2.Range("B17").Value = "Y"
$suite2.Range("B18").Value = "Y"
$suite2.Range("B19").Value = "Y"
$suite2.Range("B20").Value = "Y"

# ---------------------------------------------------------------------------
# Suite3 sheet: flip Run Flag (col B) from N to Y for rows 18-24 and 26
# (row 25 intentionally stays N)
# ---------------------------------------------------------------------------
$suite3 = $wb.Worksheets.Item("Suite3")
$suite3.Range("B18").Value = "Y"
$suite3.Range("B19").Value = "Y"
$suite3.Range("B20").Value = "Y"
$suite3.Range("B21").Value = "Y"
$suite3.Range("B22").Value = "Y"
$suite3.Range("B23").Value = "Y"
$suite3.Range("B24").Value = "Y"
$suite3.Range("B26").Value = "Y"

# ---------------------------------------------------------------------------
# Row heights that re-wrap to an extra text line because of the new content
# ---------------------------------------------------------------------------
$suite2.Rows.Item(9).RowHeight = 41.95
$suite2.Rows.Item(12).RowHeight = 41.95
$suite2.Rows.Item(13).RowHeight = 41.95
$suite2.Rows.Item(20).RowHeight = 41.95

$suite3.Rows.Item(19).RowHeight = 41.95
$suite3.Rows.Item(26).RowHeight = 41.75

# ---------------------------------------------------------------------------
# Selection / active-cell bookkeeping, matching where the editor left the
# cursor on each sheet after making the changes above.
# ---------------------------------------------------------------------------
$suite1.Range("B10").Select()
$suite2.Range("B1").Select()
$suite3.Range("B27").Select()
